# Apply "major staff excel changes":
#  - Rows 3, 6 and 7 (Sasikala / Sivaraj / Kavitha) are reshuffled:
#      new Row 3 <- old Row 7 data (Kavitha), with her title upgraded to "Dr."
#      new Row 6 <- old Row 3 data (Sasikala)
#      new Row 7 <- old Row 6 data (Sivaraj)
#  - Row 9 (Balamurugapandian) gets a title fix: "Dr BALAMURUGAPANDIAN N" -> "Dr. BALAMURUGAPANDIAN N"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Capture the "before" values of the three rows that get shuffled (columns A through J).
$cols = @(1,2,3,4,5,6,7,8,9,10)

$row3 = @{}
$row6 = @{}
$row7 = @{}
foreach ($c in $cols) {
    $row3[$c] = $ws.Cells.Item(3, $c).Value2
    $row6[$c] = $ws.Cells.Item(6, $c).Value2
    $row7[$c] = $ws.Cells.Item(7, $c).Value2
}

# New Row 3 = old Row 7 (Kavitha), title changed from "Mrs." to "Dr."
foreach ($c in $cols) {
    $ws.Cells.Item(3, $c).Value = $row7[$c]
}
$ws.Cells.Item(3, 1).Value = "Dr. KAVITHA K"

# New Row 6 = old Row 3 (Sasikala)
foreach ($c in $cols) {
    $ws.Cells.Item(6, $c).Value = $row3[$c]
}

# New Row 7 = old Row 6 (Sivaraj)
foreach ($c in $cols) {
    $ws.Cells.Item(7, $c).Value = $row6[$c]
}

# Row 9: fix the name to include the missing period after "Dr"
$ws.Cells.Item(9, 1).Value = "Dr. BALAMURUGAPANDIAN N"
